$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country stats (column order: A=Pais, B=Casos totales, C=Nuevos casos,
#     D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) ---

# Tailandia (row 117)
$ws.Cells.Item(117,2).Value2 = 3381
$ws.Cells.Item(117,3).Value2 = 3
$ws.Cells.Item(117,4).Value2 = 3198
$ws.Cells.Item(117,5).Value2 = 125

# Belice (row 170)
$ws.Cells.Item(170,2).Value2 = 475
$ws.Cells.Item(170,3).Value2 = 23
$ws.Cells.Item(170,4).Value2 = 38
$ws.Cells.Item(170,5).Value2 = 433
$ws.Cells.Item(170,7).Value2 = 1
$ws.Cells.Item(170,8).Value2 = 4

# Mongolia (row 181)
$ws.Cells.Item(181,4).Value2 = 278
$ws.Cells.Item(181,5).Value2 = 20

# Granada (row 204)
$ws.Cells.Item(204,4).Value2 = 24
$ws.Cells.Item(204,5).Value2 = 0

# Montserrat now sorts ahead of Islas Malvinas -> swap the two rows (213/214)
# in place, including the country label, so the row order/content matches.
$ws.Cells.Item(213,1).Value2 = "Montserrat"
$ws.Cells.Item(213,4).Value2 = 12
$ws.Cells.Item(213,8).Value2 = 1

$ws.Cells.Item(214,1).Value2 = "Islas Malvinas"
$ws.Cells.Item(214,4).Value2 = 13
$ws.Cells.Item(214,8).Value2 = 0

# --- Update the "last refreshed" timestamp banner ---
$ws.Cells.Item(1,1).Value2 = "Datos actualizados a 18 de Agosto de 2020 a las 06:39"
